$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leetcode")

# Fix typo / normalize existing cell F37 (LISTS column) to the canonical
# "Neetcode 150" spelling (it previously held a duplicate mis-spelled
# string "Neetcod 150").
$ws.Range("F37").Value = "Neetcode 150"

# New row 38: Reorder List
$ws.Range("A38").Value = "Leetcode"
$ws.Range("B38").Value = 143
$ws.Range("C38").Value = "Reorder List"
$ws.Range("D38").Value = "Linked Lists"
$ws.Range("E38").Value = "Medium"
$ws.Range("F38").Value = "Neetcode 150"
$ws.Range("G38").Value = "STRUGGLED"
$ws.Range("H38").Value = "18/06/2025"
$ws.Range("I38").Value = "Using an array seemed to help a lot. Trick was to remember how saving linked list in place works."

# New row 39: Remove Nth Node from End of List
$ws.Range("A39").Value = "Leetcode"
$ws.Range("B39").Value = 19
$ws.Range("C39").Value = "Remove Nth Node from End of List"
$ws.Range("D39").Value = "Linked Lists"
$ws.Range("E39").Value = "Medium"
$ws.Range("F39").Value = "Neetcode 150"

$ws.Range("A39").Select()
